# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet 1): refresh "want to go" counters -----------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 110
$ws1.Range("F3").Value  = 12154
$ws1.Range("F4").Value  = 52
$ws1.Range("F5").Value  = 237
$ws1.Range("F8").Value  = 12072
$ws1.Range("F9").Value  = 507
$ws1.Range("F11").Value = 116
$ws1.Range("F12").Value = 602
$ws1.Range("F14").Value = 5948
$ws1.Range("F15").Value = 136
$ws1.Range("F16").Value = 3563
$ws1.Range("F17").Value = 206

# --- Sheet "演出" (sheet 2): the 2024-06-01 event ended, drop its row ------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows("2:2").Delete()
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("F3").Value = 11

# --- Sheet "全部类型" (sheet 4): same row drop + refreshed counters -------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows("2:2").Delete()
for ($i = 2; $i -le 21; $i++) {
    $ws4.Range("A$i").Value = $i - 1
}
$ws4.Range("F2").Value  = 110
$ws4.Range("F4").Value  = 12154
$ws4.Range("F5").Value  = 52
$ws4.Range("F6").Value  = 237
$ws4.Range("F7").Value  = 11
$ws4.Range("F10").Value = 12072
$ws4.Range("F11").Value = 507
$ws4.Range("F13").Value = 116
$ws4.Range("F14").Value = 602
$ws4.Range("F17").Value = 5948
$ws4.Range("F18").Value = 136
$ws4.Range("F19").Value = 3563
$ws4.Range("F20").Value = 206
